# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every "protocol" worksheet (i.e. every sheet except the first five
# overview/reference sheets), prepend the worksheet's own name to the
# value in column A for each data row (row 1 is the "Name" header and is
# left untouched). e.g. on sheet "free1", "Step4 Seed" -> "free1 Step4 Seed".

$wb = $excel.ActiveWorkbook

# The first five sheets (JessicaFPJourney, NRWaves, PersonalJessicaFP,
# PositiveSpin, ReEngagement) are overview/reference sheets and are not
# touched by this change - only the protocol sheets starting at index 6
# (price1) through the last sheet (boosters) get the prefix.
$skipSheets = @("JessicaFPJourney", "NRWaves", "PersonalJessicaFP", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    if ($skipSheets -contains $ws.Name) {
        continue
    }

    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2

        if ($null -eq $current) {
            continue
        }

        $currentText = [string]$current

        if ($currentText.Length -eq 0) {
            continue
        }

        $prefix = $ws.Name + " "

        # Guard against double-prefixing if this were ever run twice.
        if (-not $currentText.StartsWith($prefix)) {
            $cell.Value = $prefix + $currentText
        }
    }
}
